# Weekly update: insert this week's two new price records for
# "Vega Monumental Concepción - Alcachofa" at the top of the data block
# (rows 23-24), pushing the existing historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 23 and 24; everything that was there
# (previously rows 23-43) shifts down to rows 25-45.
$ws.Range("A23:A24").EntireRow.Insert()

# New row 23
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = 'Vega Monumental Concepción'
$ws.Range("C23").Value = 'Bíobío'
$ws.Range("D23").Value = 44741
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112013
$ws.Range("G23").Value = 'Alcachofa'
$ws.Range("H23").Value = 'Argentina(o)'
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 16000
$ws.Range("L23").Value = 17000
$ws.Range("M23").Value = 16500
$ws.Range("N23").Value = '$/caja 50 unidades'
$ws.Range("O23").Value = 'Provincia de Limarí'
$ws.Range("P23").Value = 330
$ws.Range("Q23").Value = 50
$ws.Range("R23").Value = 'Hortaliza'

# New row 24
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = 'Vega Monumental Concepción'
$ws.Range("C24").Value = 'Bíobío'
$ws.Range("D24").Value = 44741
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 100112013
$ws.Range("G24").Value = 'Alcachofa'
$ws.Range("H24").Value = 'Española'
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 21000
$ws.Range("N24").Value = '$/caja 30 unidades'
$ws.Range("O24").Value = 'Provincia de Limarí'
$ws.Range("P24").Value = 700
$ws.Range("Q24").Value = 30
$ws.Range("R24").Value = 'Hortaliza'
